# Update "想去人数" (number of people interested) values for several
# events across the "展览", "演出" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 364
$ws1.Range("F3").Value = 774
$ws1.Range("F4").Value = 271
$ws1.Range("F5").Value = 825
$ws1.Range("F6").Value = 2034
$ws1.Range("F7").Value = 182

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 13

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 364
$ws4.Range("F3").Value = 774
$ws4.Range("F4").Value = 271
$ws4.Range("F5").Value = 13
$ws4.Range("F7").Value = 825
$ws4.Range("F8").Value = 2034
$ws4.Range("F10").Value = 182
